$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 450
$ws.Range("B4").Value = 150
$ws.Range("B5").Value = 150
$ws.Range("B6").Value = 450
